{"js": "// Change \"F values of the existing methods...\" to \"F-values of the existing\n// methods...\" (insert a hyphen between \"F\" and \"values\"), and keep the\n// document's \"_GoBack\" last-edit bookmark positioned at the edit, i.e.\n// right after the newly inserted hyphen.\n\nconst body = context.document.body;\n\n// The \"_GoBack\" bookmark marks the location of the most recent edit.\n// Drop the existing one (wherever it currently sits) before we make our\n// edit; we'll re-add it at the new edit location below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Locate the sentence fragment that needs the hyphen inserted. Using the\n// full phrase keeps the match unambiguous (there's an unrelated\n// \"F-values\" occurrence elsewhere in the document, in a figure caption).\nconst target = body.search(\n  \"F values of the existing methods and the proposed method\",\n  { matchCase: true }\n);\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\n    \"F-values of the existing methods and the proposed method\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Re-insert the \"_GoBack\" bookmark right after \"F-\" (before \"values\"),\n// matching where Word leaves it after the edit.\nconst anchor = body.search(\"comparison between the F-\", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const insertionPoint = anchor.items[0].getRange(\"End\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Change \"F values of the existing methods...\" to \"F-values of the existing\n# methods...\" (insert a hyphen between \"F\" and \"values\"), and keep the\n# document's \"_GoBack\" last-edit bookmark positioned at the edit, i.e.\n# right after the newly inserted hyphen.\n\n$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark marks the location of the most recent edit. Drop\n# the existing one (wherever it currently sits) before we make our edit;\n# we'll re-add it at the new edit location below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Find the sentence fragment and replace \"F values\" with \"F-values\". Using\n# the full phrase keeps the match unambiguous (there's an unrelated\n# \"F-values\" occurrence elsewhere in the document, in a figure caption).\n$find = $d.Content.Find\n$find.Text = \"F values of the existing methods and the proposed method\"\n$find.Replacement.Text = \"F-values of the existing methods and the proposed method\"\n$find.Execute([Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceOne) | Out-Null\n\n# Re-create the \"_GoBack\" bookmark right after \"F-\" (before \"values\"),\n# matching where Word leaves it after the edit.\n$rng = $d.Content.Duplicate\n$find2 = $rng.Find\n$find2.Text = \"comparison between the F-\"\nif ($find2.Execute()) {\n    $rng.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n}\n"}
